# "deactivated empty shipping instructions card"
# Clears the sample/demo values that had been filled into the Shipping
# Instruction template, leaving the field labels and formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BillOfLading")

# --- Consignor (Shipper) details block ---
$ws.Range("C4:C9").ClearContents()

# --- Consignee / Notify Party details block ---
$ws.Range("C11:C16").ClearContents()
$ws.Range("F11:F16").ClearContents()

# --- Container / cargo table header row 18 stays, data rows cleared ---
$ws.Range("B19:H19").ClearContents()
$ws.Range("B20:D27").ClearContents()

# Reflect the cell the author last had selected after clearing the data
$ws.Range("D19").Select()
